$wb = $excel.ActiveWorkbook

# --- Sheet: Marking Scheme ---
$wsMarking = $wb.Worksheets.Item("Marking Scheme")
$wsMarking.Range("B2").Value = 'The VTC is the largest provider of VPET in Hong Kong. Briefly explain what VPET stands for and why it is important for Hong Kong’s workforce development.'
$wsMarking.Range("B3").Value = 'Compare IVE (Hong Kong Institute of Vocational Education) and THEi (Technological and Higher Education Institute of Hong Kong). What is the main difference between the types of qualifications/programmes offered by these two institutions?'
$wsMarking.Range("B4").Value = 'VTC emphasizes the " Think and Do" approach. Explain what this phrase means in the context of a student''s learning experience.'
$wsMarking.Range("B5").Value = 'If a Secondary 6 student does not achieve the minimum entrance requirements for a Bachelor''s Degree or a Higher Diploma, what is the VTC study pathway available to them to eventually reach a Higher Diploma level? (Name the specific foundation programme).'
$wsMarking.Range("B6").Value = 'Why does the VTC collaborate closely with industry partners (companies and trade associations)? Give two examples of how this benefits students.'
$wsMarking.Range("C2").Value = '- **Vocational and Professional Education and Training** (2 marks)
- Focus on **practical skills** or **specialized trades** (4 marks)
- Benefit to workforce: reducing skills gap, employment readiness (4 marks)

**General Grading Principles:**
- **9-10 marks**: Complete, accurate, correct terminology.
- **6-8 marks**: Mostly correct, misses detail.
- **3-5 marks**: Basic understanding only.
- **0-2 marks**: Incorrect or irrelevant.

---

**General Grading Guide:**
- **9-10 marks**: The answer is complete, accurate, uses correct terminology, and is well-explained.
- **6-8 marks**: The answer is mostly correct but misses a specific detail (e.g., forgets the full name of a diploma) or the explanation is slightly vague.
- **3-5 marks**: The student shows basic understanding but misses the core point or only answers half the question.
- **0-2 marks**: The answer is largely incorrect, irrelevant, or blank.'
$wsMarking.Range("C3").Value = '- **IVE**: Primarily focuses on **Higher Diploma (HD)** programmes which are practical and technical in nature (5 marks)
- **THEi**: Focuses on vocationally-oriented **Bachelor’s Degree** programmes that combine practical application with higher-level theory (5 marks)

**General Grading Principles:**
- **9-10 marks**: Complete, accurate, correct terminology.
- **6-8 marks**: Mostly correct, misses detail.
- **3-5 marks**: Basic understanding only.
- **0-2 marks**: Incorrect or irrelevant.

---

**General Grading Guide:**
- **9-10 marks**: The answer is complete, accurate, uses correct terminology, and is well-explained.
- **6-8 marks**: The answer is mostly correct but misses a specific detail (e.g., forgets the full name of a diploma) or the explanation is slightly vague.
- **3-5 marks**: The student shows basic understanding but misses the core point or only answers half the question.
- **0-2 marks**: The answer is largely incorrect, irrelevant, or blank.'
$wsMarking.Range("C4").Value = '- **Think**: Theory/Academic knowledge/Brainpower (3 marks)
- **Do**: Practical skills/Hands-on/Technical execution (3 marks)
- **Synthesis**: Ability to solve problems by combining both head and hands (4 marks)

**General Grading Principles:**
- **9-10 marks**: Complete, accurate, correct terminology.
- **6-8 marks**: Mostly correct, misses detail.
- **3-5 marks**: Basic understanding only.
- **0-2 marks**: Incorrect or irrelevant.

---

**General Grading Guide:**
- **9-10 marks**: The answer is complete, accurate, uses correct terminology, and is well-explained.
- **6-8 marks**: The answer is mostly correct but misses a specific detail (e.g., forgets the full name of a diploma) or the explanation is slightly vague.
- **3-5 marks**: The student shows basic understanding but misses the core point or only answers half the question.
- **0-2 marks**: The answer is largely incorrect, irrelevant, or blank.'
$wsMarking.Range("C5").Value = '- Correctly naming the **Diploma of Foundation Studies (DFS)** or **Diploma of Vocational Education (DVE)** (5 marks)
- Explanation of progression: Successful completion allows entry to **Higher Diploma (HD)** programmes (5 marks)

**General Grading Principles:**
- **9-10 marks**: Complete, accurate, correct terminology.
- **6-8 marks**: Mostly correct, misses detail.
- **3-5 marks**: Basic understanding only.
- **0-2 marks**: Incorrect or irrelevant.

---

**General Grading Guide:**
- **9-10 marks**: The answer is complete, accurate, uses correct terminology, and is well-explained.
- **6-8 marks**: The answer is mostly correct but misses a specific detail (e.g., forgets the full name of a diploma) or the explanation is slightly vague.
- **3-5 marks**: The student shows basic understanding but misses the core point or only answers half the question.
- **0-2 marks**: The answer is largely incorrect, irrelevant, or blank.'
$wsMarking.Range("C6").Value = '- **General Explanation**: Ensures curriculum is up-to-date with market trends and industry needs (4 marks)
- **Benefit 1**: e.g., **Internship** or Work-integrated learning opportunities (3 marks)
- **Benefit 2**: e.g., **Job placement support** or access to industry-standard equipment (3 marks)

**General Grading Principles:**
- **9-10 marks**: Complete, accurate, correct terminology.
- **6-8 marks**: Mostly correct, misses detail.
- **3-5 marks**: Basic understanding only.
- **0-2 marks**: Incorrect or irrelevant.

---

**General Grading Guide:**
- **9-10 marks**: The answer is complete, accurate, uses correct terminology, and is well-explained.
- **6-8 marks**: The answer is mostly correct but misses a specific detail (e.g., forgets the full name of a diploma) or the explanation is slightly vague.
- **3-5 marks**: The student shows basic understanding but misses the core point or only answers half the question.
- **0-2 marks**: The answer is largely incorrect, irrelevant, or blank.'

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B7").Value = '2026-01-07 02:18:41'

# --- Sheet: Question Overview ---
$wsOverview = $wb.Worksheets.Item("Question Overview")
$wsOverview.Range("B2").Value = 'The VTC is the largest provider of VPET in Hong Kong. Briefly explain what VPET stands for and why it is important for Hong Kong’s workforce development.'
$wsOverview.Range("B3").Value = 'Compare IVE (Hong Kong Institute of Vocational Education) and THEi (Technological and Higher Education Institute of Hong Kong). What is the main difference between the types of qualifications/programmes offered by these two institutions?'
$wsOverview.Range("B4").Value = 'VTC emphasizes the " Think and Do" approach. Explain what this phrase means in the context of a student''s learning experience.'
$wsOverview.Range("B5").Value = 'If a Secondary 6 student does not achieve the minimum entrance requirements for a Bachelor''s Degree or a Higher Diploma, what is the VTC study pathway available to them to eventually reach a Higher Diploma level? (Name the specific foundation programme).'
$wsOverview.Range("B6").Value = 'Why does the VTC collaborate closely with industry partners (companies and trade associations)? Give two examples of how this benefits students.'
